$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the account numbers in column A (rows 2-6) to the new sequence.
# These look like numbers but are stored as shared strings (text) in the
# workbook, so we temporarily force Text formatting while assigning the
# values (otherwise Excel would auto-convert them to numeric cells), then
# clear the formatting again so the cells end up without any explicit
# style, exactly as they were before the edit.
$rng = $ws.Range("A2:A6")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "1008617924"
$ws.Range("A3").Value = "1008617925"
$ws.Range("A4").Value = "1008617926"
$ws.Range("A5").Value = "1008617927"
$ws.Range("A6").Value = "1008617928"

$rng.ClearFormats()
